# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" worksheets, reflecting refreshed figures.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1474
$ws1.Range("F4").Value  = 1761
$ws1.Range("F6").Value  = 145
$ws1.Range("F7").Value  = 657
$ws1.Range("F8").Value  = 34
$ws1.Range("F9").Value  = 63
$ws1.Range("F12").Value = 80
$ws1.Range("F13").Value = 149
$ws1.Range("F18").Value = 4831
$ws1.Range("F19").Value = 45
$ws1.Range("F21").Value = 106
$ws1.Range("F22").Value = 2215
$ws1.Range("F25").Value = 2072

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1474
$ws4.Range("F4").Value  = 1761
$ws4.Range("F6").Value  = 145
$ws4.Range("F7").Value  = 657
$ws4.Range("F8").Value  = 34
$ws4.Range("F9").Value  = 63
$ws4.Range("F12").Value = 80
$ws4.Range("F13").Value = 149
$ws4.Range("F18").Value = 4831
$ws4.Range("F20").Value = 45
$ws4.Range("F23").Value = 106
$ws4.Range("F24").Value = 2215
$ws4.Range("F27").Value = 2072
